$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F334").Value = 196209
$ws.Range("F335").Value = 131122
$ws.Range("G335").Value = 3006
$ws.Range("F336").Value = 101986
$ws.Range("F337").Value = 104155
$ws.Range("F338").Value = 227269
$ws.Range("F339").Value = 660826
$ws.Range("G339").Value = 5500
$ws.Range("F340").Value = 385549
$ws.Range("G340").Value = 3311
$ws.Range("F341").Value = 291467
$ws.Range("G341").Value = 3672
$ws.Range("F342").Value = 179655
$ws.Range("G342").Value = 3073
$ws.Range("F343").Value = 132649
$ws.Range("F344").Value = 135670
$ws.Range("F345").Value = 292098
$ws.Range("G345").Value = 3325
$ws.Range("F346").Value = 672802
$ws.Range("G346").Value = 4793
$ws.Range("F347").Value = 343769
$ws.Range("G347").Value = 2918
$ws.Range("F348").Value = 232320
$ws.Range("F350").Value = 127074
$ws.Range("G350").Value = 2781
$ws.Range("F351").Value = 150645
$ws.Range("G352").Value = 3543
$ws.Range("F353").Value = 724028
$ws.Range("G353").Value = 5267
$ws.Range("F354").Value = 311401
$ws.Range("G354").Value = 2849
$ws.Range("F355").Value = 222231
$ws.Range("G355").Value = 3457
$ws.Range("F356").Value = 160293
$ws.Range("G356").Value = 2887
$ws.Range("F357").Value = 138350
$ws.Range("G357").Value = 3027
$ws.Range("F358").Value = 157293
$ws.Range("F359").Value = 321139
$ws.Range("G359").Value = 3350
$ws.Range("F360").Value = 747214
$ws.Range("G360").Value = 5112
$ws.Range("F361").Value = 332591
$ws.Range("G361").Value = 2622
$ws.Range("F362").Value = 227832
$ws.Range("F363").Value = 188477
$ws.Range("G363").Value = 2768
$ws.Range("F364").Value = 167652
$ws.Range("G364").Value = 2470
$ws.Range("F365").Value = 181231
$ws.Range("G365").Value = 2365
$ws.Range("F366").Value = 338313
$ws.Range("G366").Value = 2839
$ws.Range("F367").Value = 763909
$ws.Range("G367").Value = 3900
$ws.Range("F368").Value = 344419
$ws.Range("G368").Value = 2290
$ws.Range("F369").Value = 235155
$ws.Range("G369").Value = 2596
$ws.Range("F370").Value = 181799
$ws.Range("G370").Value = 2028
$ws.Range("F371").Value = 158484
$ws.Range("G371").Value = 1948
$ws.Range("F372").Value = 177105
$ws.Range("G372").Value = 1842
$ws.Range("F373").Value = 345196
$ws.Range("G373").Value = 2351
$ws.Range("F374").Value = 768503
$ws.Range("G374").Value = 3402
$ws.Range("F375").Value = 350514
$ws.Range("G375").Value = 1844
$ws.Range("F376").Value = 219706
$ws.Range("G376").Value = 2202
$ws.Range("F377").Value = 175272
$ws.Range("G377").Value = 1818
$ws.Range("F378").Value = 155507
$ws.Range("G378").Value = 1521
$ws.Range("F379").Value = 177240
$ws.Range("G379").Value = 1600
$ws.Range("F380").Value = 339265
$ws.Range("G380").Value = 1973
$ws.Range("F381").Value = 729476
$ws.Range("G381").Value = 2633
$ws.Range("F382").Value = 352453
$ws.Range("G382").Value = 1600
$ws.Range("F383").Value = 216184
$ws.Range("G383").Value = 1727
$ws.Range("F384").Value = 164489
$ws.Range("G384").Value = 1475